$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.565.60"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.451.37"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.76"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.68"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.450.80"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.63"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.391"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.030.85"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.68"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.31%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.444.03"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.591.75"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.41"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "400.72"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.566"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.583.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.66"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.49"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -6.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.25"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  +0.60%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "24.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.07"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.471.09"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.35"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0792"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "28.15"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.803"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.55"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.12%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.74"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.633.57"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.16"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.42%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.97"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.14"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.41"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.93%  "
